# Apply the "Final code worked as team" commit to the Allergies sheet:
#   - rename the sheet from "Filter -1 Allergies - Bonus Poi" to "Final List for Allergies"
#   - retitle the header cell from "Allergies (Bonus points)" to "Allergies"
#   - insert a new "Eliminate" header row (row 2), matching the header row used
#     on the other two sheets, shifting the allergen list down by one row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Filter -1 Allergies - Bonus Poi")

# Rename the worksheet tab.
$ws.Name = "Final List for Allergies"

# Update the sheet title in A1.
$ws.Range("A1").Value = "Allergies"

# Insert a new row right below the title so the allergen rows move down one.
$ws.Rows.Item(2).Insert()

# Give the new row the same look (bold/peach/bordered) as the existing title row.
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row with the "Eliminate" header, like the other two sheets.
$ws.Range("A2").Value = "Eliminate"
